# Update the three-digit-number x one-digit-number practice problems
# to the newly generated set of operands.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "512×2="; New = "218×5=" },
    @{ Old = "385×3="; New = "995×8=" },
    @{ Old = "869×5="; New = "354×2=" },
    @{ Old = "388×2="; New = "762×3=" },
    @{ Old = "702×9="; New = "417×7=" },
    @{ Old = "603×3="; New = "669×6=" },
    @{ Old = "187×3="; New = "294×6=" },
    @{ Old = "390×2="; New = "375×9=" },
    @{ Old = "714×4="; New = "691×4=" },
    @{ Old = "716×7="; New = "169×5=" },
    @{ Old = "359×2="; New = "963×9=" },
    @{ Old = "203×7="; New = "829×4=" },
    @{ Old = "898×6="; New = "718×5=" },
    @{ Old = "591×8="; New = "381×3=" },
    @{ Old = "871×9="; New = "605×6=" },
    @{ Old = "768×2="; New = "185×3=" },
    @{ Old = "108×8="; New = "672×8=" },
    @{ Old = "947×4="; New = "438×9=" },
    @{ Old = "574×9="; New = "660×6=" },
    @{ Old = "503×9="; New = "251×6=" },
    @{ Old = "875×6="; New = "186×8=" },
    @{ Old = "617×2="; New = "379×5=" },
    @{ Old = "985×9="; New = "354×3=" },
    @{ Old = "560×8="; New = "912×5=" },
    @{ Old = "640×8="; New = "106×9=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2) | Out-Null
}
